# A new review/warranty-evaluation entry was inserted at row 11, pushing the
# existing rows 11-21 down to rows 12-22 (dimension grows from D21 to D22).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 11, shifting rows 11..21 down to 12..22.
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with the new record's data.
$ws.Cells.Item(11, 1).Value = 5
$ws.Cells.Item(11, 2).Value = ""
$ws.Cells.Item(11, 3).Value = 45951.66084622685
$ws.Cells.Item(11, 4).Value = "MjRkNTkzODMtN2IzMC00N2JhLWI1ZDQtNjYwNDFhNjUxZTU1OjU3MDE2"
